$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Berenjena, Vega Monumental Concepción) is
# inserted as row 21; every existing record from row 21 downwards shifts
# down by one row, and the sheet's used range grows from R87 to R88.
$ws.Rows.Item(21).Insert()

$ws.Cells.Item(21, 1).Value = 11
$ws.Cells.Item(21, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(21, 3).Value = "Bíobío"
$ws.Cells.Item(21, 4).Value = 44701
$ws.Cells.Item(21, 5).Value = 8
$ws.Cells.Item(21, 6).Value = 100112001
$ws.Cells.Item(21, 7).Value = "Berenjena"
$ws.Cells.Item(21, 8).Value = "Sin especificar"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 150
$ws.Cells.Item(21, 11).Value = 6000
$ws.Cells.Item(21, 12).Value = 6500
$ws.Cells.Item(21, 13).Value = 6233
$ws.Cells.Item(21, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(21, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(21, 16).Value = 104
$ws.Cells.Item(21, 17).Value = 60
$ws.Cells.Item(21, 18).Value = "Hortaliza"
